$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D2: change from shared string "dsdsd" to numeric 2250
$ws.Range("D2").Value = 2250

# Update B3: add boolean TRUE
$ws.Range("B3").Value = $true

# Update D3: change from 2222 to 1
$ws.Range("D3").Value = 1

# Update selection to B3
$ws.Activate()
$ws.Range("B3").Select()
